# feat: add 2022-Q1 data
#
# 1. Insert a new worksheet "2022-Q1" positioned after "2021-Q4" and
#    before "总计", containing the fund holding detail for the new quarter.
# 2. Update the "总计" (totals) worksheet by inserting a new top data row
#    for "2022-Q1" and shifting the existing rows (and their running index
#    in column A) down by one.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: create the "2022-Q1" sheet
# ---------------------------------------------------------------------
# Duplicate the "2021-Q4" sheet (same column layout/formatting) right
# after itself -- this keeps sheetPr/pageMargins/styles consistent with
# the other per-quarter sheets -- then rename it and replace its data.
$q4Sheet = $wb.Worksheets.Item("2021-Q4")
$q4Sheet.Copy($null, $q4Sheet)
$newSheet = $wb.Worksheets.Item("2021-Q4 (2)")
$newSheet.Name = "2022-Q1"

# Row 2: 004209 / 大成智惠量化多策略灵活配置混合
$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "'004209"
$newSheet.Range("C2").Value = "大成智惠量化多策略灵活配置混合"
$newSheet.Range("D2").Value = "'1.26"
$newSheet.Range("E2").Value = "'94.24"
$newSheet.Range("F2").Value = "'6.62"
$newSheet.Range("G2").Value = "'0.0834"
$newSheet.Range("H2").Value = 5

# Row 3: 161224 / 国投瑞银新丝路灵活配置混合(LOF)
$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = "'161224"
$newSheet.Range("C3").Value = "国投瑞银新丝路灵活配置混合(LOF)"
$newSheet.Range("D3").Value = "'0.77"
$newSheet.Range("E3").Value = "'94.48"
$newSheet.Range("F3").Value = "'8.63"
$newSheet.Range("G3").Value = "'0.0665"
$newSheet.Range("H3").Value = 1

# ---------------------------------------------------------------------
# Step 2: update the "总计" sheet with the new 2022-Q1 summary row
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")

# Push the existing data rows (old rows 2-6) down by one row.
$totalSheet.Range("A2:D2").Insert()

# New row 2: 2022-Q1 summary.
$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 2
$totalSheet.Range("D2").Value = 0.15

# Re-number the running index in column A for the rows that shifted down.
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("A5").Value = 3
$totalSheet.Range("A6").Value = 4
$totalSheet.Range("A7").Value = 5

# Match the formatting of the other rows: column A keeps the index style,
# columns B:D use the plain/default style (no quote-prefix style leakage
# from the row-insert operation).
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)
$totalSheet.Range("B2:D2").ClearFormats()

Write-Host "2022-Q1 sheet added and 总计 sheet updated"
